$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 15049
$ws.Range("E2").Value = 1443
$ws.Range("F2").Value = 1443
$ws.Range("G2").Value = 1788
$ws.Range("H2").Value = 1447
$ws.Range("I2").Value = 1447
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 117838
$ws.Range("L2").Value = 107066
$ws.Range("M2").Value = 10771
$ws.Range("N2").Value = 8998
$ws.Range("O2").Value = 1774
$ws.Range("P2").Value = 3118
$ws.Range("Q2").Value = -13077
$ws.Range("R2").Value = -3879
$ws.Range("S2").Value = 19191
$ws.Range("T2").Value = 50
$ws.Range("V2").Value = 17704
$ws.Range("W2").Value = 9.59
$ws.Range("X2").Value = 9.619999999999999
$ws.Range("Y2").Value = 17.93
$ws.Range("Z2").Value = 1.41
$ws.Range("AA2").Value = 994
$ws.Range("AB2").Value = 248.59
$ws.Range("AC2").Value = 413
$ws.Range("AD2").Value = 8.73
$ws.Range("AE2").Value = 2605
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 129
$ws.Range("AH2").Value = 3.57
$ws.Range("AI2").Value = 30.62
$ws.Range("AJ2").Value = 351549340
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 32410
$ws.Range("E3").Value = 4051
$ws.Range("F3").Value = 4051
$ws.Range("G3").Value = 4065
$ws.Range("H3").Value = 2873
$ws.Range("I3").Value = 2961
$ws.Range("J3").Value = -88
$ws.Range("K3").Value = 144287
$ws.Range("L3").Value = 127101
$ws.Range("M3").Value = 17186
$ws.Range("N3").Value = 17186
$ws.Range("P3").Value = 4966
$ws.Range("Q3").Value = -22628
$ws.Range("R3").Value = -625
$ws.Range("S3").Value = 24396
$ws.Range("T3").Value = 53
$ws.Range("V3").Value = 20401
$ws.Range("W3").Value = 12.5
$ws.Range("X3").Value = 8.869999999999999
$ws.Range("Y3").Value = 22.62
$ws.Range("Z3").Value = 2.26
$ws.Range("AA3").Value = 739.58
$ws.Range("AB3").Value = 248
$ws.Range("AC3").Value = 709
$ws.Range("AD3").Value = 5.63
$ws.Range("AE3").Value = 3790
$ws.Range("AF3").Value = 1.05
$ws.Range("AG3").Value = 230
$ws.Range("AH3").Value = 5.76
$ws.Range("AI3").Value = 35.22
$ws.Range("AJ3").Value = 496638534
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 49466
$ws.Range("E4").Value = 3269
$ws.Range("F4").Value = 3269
$ws.Range("G4").Value = 3326
$ws.Range("H4").Value = 2538
$ws.Range("I4").Value = 2538
$ws.Range("K4").Value = 172678
$ws.Range("L4").Value = 153895
$ws.Range("M4").Value = 18783
$ws.Range("N4").Value = 18783
$ws.Range("P4").Value = 4966
$ws.Range("Q4").Value = -3989
$ws.Range("R4").Value = -9962
$ws.Range("S4").Value = 15052
$ws.Range("T4").Value = 38
$ws.Range("V4").Value = 36229
$ws.Range("W4").Value = 6.61
$ws.Range("X4").Value = 5.13
$ws.Range("Y4").Value = 14.12
$ws.Range("Z4").Value = 1.6
$ws.Range("AA4").Value = 819.33
$ws.Range("AB4").Value = 280.16
$ws.Range("AC4").Value = 511
$ws.Range("AD4").Value = 6.77
$ws.Range("AE4").Value = 4142
$ws.Range("AF4").Value = 0.84
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 5.78
$ws.Range("AI4").Value = 35.73
$ws.Range("AJ4").Value = 496638534
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 52976
$ws.Range("E5").Value = 4436
$ws.Range("F5").Value = 4436
$ws.Range("G5").Value = 4664
$ws.Range("H5").Value = 3552
$ws.Range("I5").Value = 3510
$ws.Range("J5").Value = 42
$ws.Range("K5").Value = 234506
$ws.Range("L5").Value = 201380
$ws.Range("M5").Value = 33126
$ws.Range("N5").Value = 32129
$ws.Range("O5").Value = 997
$ws.Range("P5").Value = 7175
$ws.Range("Q5").Value = -36507
$ws.Range("R5").Value = 7361
$ws.Range("S5").Value = 32848
$ws.Range("T5").Value = 48
$ws.Range("V5").Value = 72487
$ws.Range("W5").Value = 8.369999999999999
$ws.Range("X5").Value = 6.71
$ws.Range("Y5").Value = 13.79
$ws.Range("Z5").Value = 1.72
$ws.Range("AA5").Value = 607.92
$ws.Range("AB5").Value = 363.05
$ws.Range("AC5").Value = 559
$ws.Range("AD5").Value = 8.18
$ws.Range("AE5").Value = 4764
$ws.Range("AF5").Value = 0.96
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 4.37
$ws.Range("AI5").Value = 36.69
$ws.Range("AJ5").Value = 605641072
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 87394
$ws.Range("E6").Value = 5323
$ws.Range("F6").Value = 5323
$ws.Range("G6").Value = 5892
$ws.Range("H6").Value = 4338
$ws.Range("I6").Value = 4282
$ws.Range("K6").Value = 311618
$ws.Range("L6").Value = 276887
$ws.Range("M6").Value = 34731
$ws.Range("N6").Value = 33734
$ws.Range("P6").Value = 7175
$ws.Range("Q6").Value = -50427
$ws.Range("R6").Value = -1394
$ws.Range("S6").Value = 51373
$ws.Range("T6").Value = 29
$ws.Range("V6").Value = 114158
$ws.Range("W6").Value = 6.09
$ws.Range("X6").Value = 4.96
$ws.Range("Y6").Value = 13
$ws.Range("Z6").Value = 1.57
$ws.Range("AA6").Value = 797.23
$ws.Range("AB6").Value = 385.42
$ws.Range("AC6").Value = 602
$ws.Range("AD6").Value = 7.11
$ws.Range("AE6").Value = 5107
$ws.Range("AF6").Value = 0.84
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 4.67
$ws.Range("AI6").Value = 32.54
$ws.Range("AJ6").Value = 605641072
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 113973
$ws.Range("E7").Value = 5837
$ws.Range("G7").Value = 6074
$ws.Range("H7").Value = 4931
$ws.Range("I7").Value = 4885
$ws.Range("K7").Value = 381360
$ws.Range("L7").Value = 343210
$ws.Range("M7").Value = 38152
$ws.Range("N7").Value = 37395
$ws.Range("P7").Value = 7175
$ws.Range("W7").Value = 5.12
$ws.Range("X7").Value = 4.33
$ws.Range("Y7").Value = 13.74
$ws.Range("Z7").Value = 1.42
$ws.Range("AA7").Value = 899.59
$ws.Range("AC7").Value = 694
$ws.Range("AD7").Value = 5.08
$ws.Range("AE7").Value = 5661
$ws.Range("AF7").Value = 0.62
$ws.Range("AG7").Value = 218
$ws.Range("AH7").Value = 6.18
$ws.Range("AI7").Value = 27.55
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 81776
$ws.Range("E8").Value = 5900
$ws.Range("G8").Value = 6531
$ws.Range("H8").Value = 4830
$ws.Range("I8").Value = 4736
$ws.Range("K8").Value = 408199
$ws.Range("L8").Value = 365865
$ws.Range("M8").Value = 42335
$ws.Range("N8").Value = 41805
$ws.Range("P8").Value = 7174
$ws.Range("W8").Value = 7.21
$ws.Range("X8").Value = 5.91
$ws.Range("Y8").Value = 11.96
$ws.Range("Z8").Value = 1.22
$ws.Range("AA8").Value = 864.21
$ws.Range("AC8").Value = 673
$ws.Range("AD8").Value = 5.44
$ws.Range("AE8").Value = 6329
$ws.Range("AF8").Value = 0.58
$ws.Range("AG8").Value = 218
$ws.Range("AH8").Value = 5.97
$ws.Range("AI8").Value = 28.46
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("E9").Value = 6075
$ws.Range("G9").Value = 6635
$ws.Range("H9").Value = 4828
$ws.Range("I9").Value = 4815
$ws.Range("K9").Value = 459065
$ws.Range("L9").Value = 413240
$ws.Range("M9").Value = 45822
$ws.Range("N9").Value = 47270
$ws.Range("P9").Value = 7170
$ws.Range("Y9").Value = 14.07
$ws.Range("Z9").Value = 1.11
$ws.Range("AA9").Value = 901.83
$ws.Range("AC9").Value = 684
$ws.Range("AD9").Value = 5.35
$ws.Range("AE9").Value = 7156
$ws.Range("AF9").Value = 0.51
$ws.Range("AG9").Value = 222
$ws.Range("AH9").Value = 6.08
$ws.Range("AI9").Value = 28.53
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
